# Append: 2025-09-14 18:22 JST
# A new listing is prepended to the top of the data table (row 2),
# the two previously-existing listings shift down by one row, and a
# brand new listing is appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks/relationships up front; row contents are
# about to be fully rewritten and new hyperlinks will be (re)created below
# in the correct final row order.
$ws.Hyperlinks.Delete()

$timestamp = "2025-09-14 18:22:08"

# --- Row 2: brand new listing -------------------------------------------
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【急募】メモリデータ管理ツール開発のプロフェッショナル募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5393508"
$ws.Range("G2").Value = 158
$ws.Range("H2").Value = "◆ツール,開発 ◇管理"

# --- Row 3: previously row 2, timestamp refreshed -----------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【業務委託】アプリ開発の継続的パートナ募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5393175"
$ws.Range("G3").Value = 93
$ws.Range("H3").Value = "◆開発 ◇アプリ"

# --- Row 4: previously row 3, timestamp refreshed -----------------------
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【医療関連】会員制サイト構築のパートナーを探しています"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5393406"
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = "◇サイト"

# --- Row 5: brand new listing appended at the bottom (no skill summary) -
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【急募】Discordボット設定とサブスク化の専門家を探しています!"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5393471"
$ws.Range("G5").Value = 18

# Re-create the hyperlinks, in row order, so relationship ids come out as
# rId1..rId4 matching F2..F5, then restore the shared "Hyperlink" cell
# style (xf already present in the workbook) on each URL cell.
$ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Value())
$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value())
$ws.Hyperlinks.Add($ws.Range("F4"), $ws.Range("F4").Value())
$ws.Hyperlinks.Add($ws.Range("F5"), $ws.Range("F5").Value())

$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"

# --- Column width tweaks -------------------------------------------------
# ColumnWidth goes through Excel's character->pixel->character round trip,
# which adds a fractional offset unless the input already lands on a
# pixel boundary; subtracting ~0.86 lands cleanly back on the integer
# widths used by the target workbook (B: 29->36, H: 12->13).
$ws.Columns.Item(2).ColumnWidth = 36 - 0.86
$ws.Columns.Item(8).ColumnWidth = 13 - 0.86
